$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell E6 ("It shows Company details ...") gets an updated list of buttons:
# the old text described a single "Save changes" action; the new text
# reflects that Create Supplier now offers two distinct save actions
# ("Save supplier only" and "Save supplier and Create new user").
$prefix = 'It shows Company details '
$bold = '"Checkbox of Active, Company market dropdown, Company name, Also known as, Company Reg.No, Address, Logo, Supplier email, Supplier Phone, Short description, Long description, Company time zone, Subscription information, Send subscription related notifications to, Pricing, Inventory, Email upcoming delivery schedule to, Payments, Notifications, Activity reports, Cancel , Save supplier only and Save supplier and Create new user"'
$newText = $prefix + $bold

$rng = $ws.Range("E6")
$rng.Value = $newText

# Re-apply the rich-text formatting: the leading phrase stays regular,
# the quoted list of controls is bold (matches the rest of the sheet's style).
$boldChars = $rng.Characters($prefix.Length + 1, $bold.Length)
$boldChars.Font.Bold = $true
$boldChars.Font.Name = "Calibri"
$boldChars.Font.Size = 11
$boldChars.Font.Color = 0

# Reflect the cell the author last edited/selected.
[void]$rng.Select()
